# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2410"
#   "<header>_new" -> "<header>_FV2504"
# Then wrap the data range in an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1, columns A:J = "_old" -> "_FV2410",
#        columns L:U = "_new" -> "_FV2504"; column K ("diff") is unchanged).
$oldSuffix = "_old"
$newSuffixLeft = "_FV2410"
$newSuffix = "_new"
$newSuffixRight = "_FV2504"

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    if ($null -eq $text) { continue }
    if ($text -like "*$oldSuffix") {
        $base = $text.Substring(0, $text.Length - $oldSuffix.Length)
        $cell.Value = "$base$newSuffixLeft"
    } elseif ($text -like "*$newSuffix") {
        $base = $text.Substring(0, $text.Length - $newSuffix.Length)
        $cell.Value = "$base$newSuffixRight"
    }
}

# --- 2. Turn the used range into a native Excel Table ("Table1").
$dataRange = $ws.Range("A1:U74")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
